$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update USERID column (G) from "Putri" to numeric userid 32382
$ws.Range("G2").Value = 32382
$ws.Range("G3").Value = 32382
$ws.Range("G4").Value = 32382
$ws.Range("G5").Value = 32382

# Update PASSWORD column (H) from "bni1234/" to "bni1234"
$ws.Range("H2").Value = "bni1234"
$ws.Range("H3").Value = "bni1234"
$ws.Range("H4").Value = "bni1234"
$ws.Range("H5").Value = "bni1234"

# Update PREPARATION column (F) text to reflect new userid/password
$ws.Range("F2").Value = "Username : 32382;`nPassword : bni1234;`nKode Status Peserta : 16;`nNama Status Peserta : Pensiun Lebih Awal"
$ws.Range("F3").Value = "Username : 32382;`nPassword : bni1234;`nKode Status Peserta : 16"
$ws.Range("F4").Value = "Username : 32382;`nPassword : bni1234;`nKode Status Peserta : 16;`nNama Status Peserta : Pensiun Dini"
$ws.Range("F5").Value = "Username : 32382;`nPassword : bni1234;`nKode Status Peserta : 16"

# Update the active selection to F6
$ws.Range("F6").Select()
